# calorimetry : scripts : tests : updated
#
# The "metrics" sheet gets its NRMSE value recomputed and a new RMSE
# metric row appended below SMAPE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metrics")

# NRMSE (row 3) recomputed value.
$ws.Range("B3").Value = 0.0143489666385373

# New RMSE row appended after the existing SMAPE row (row 4).
$ws.Range("A5").Value = "RMSE"
$ws.Range("B5").Value = 0.233957389972189
